# Auto-generated Excel COM-interop script to append sensor log rows
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A80:A96").NumberFormat = "@"

$ws.Cells.Item(80, 1).Value = "2026-02-04"
$ws.Cells.Item(80, 2).Value = "14:07:51"
$ws.Cells.Item(80, 3).Value = "14:00"
$ws.Cells.Item(80, 4).Value = "Bathroom"
$ws.Cells.Item(80, 5).Value = "No Motion"
$ws.Cells.Item(80, 6).Value = "Inactive"

$ws.Cells.Item(81, 1).Value = "2026-02-04"
$ws.Cells.Item(81, 2).Value = "14:07:52"
$ws.Cells.Item(81, 3).Value = "14:00"
$ws.Cells.Item(81, 4).Value = "Bathroom"
$ws.Cells.Item(81, 5).Value = "No Motion"
$ws.Cells.Item(81, 6).Value = "Inactive"

$ws.Cells.Item(82, 1).Value = "2026-02-04"
$ws.Cells.Item(82, 2).Value = "14:07:53"
$ws.Cells.Item(82, 3).Value = "14:00"
$ws.Cells.Item(82, 4).Value = "Bathroom"
$ws.Cells.Item(82, 5).Value = "No Motion"
$ws.Cells.Item(82, 6).Value = "Inactive"

$ws.Cells.Item(83, 1).Value = "2026-02-04"
$ws.Cells.Item(83, 2).Value = "14:07:53"
$ws.Cells.Item(83, 3).Value = "14:00"
$ws.Cells.Item(83, 4).Value = "Bathroom"
$ws.Cells.Item(83, 5).Value = "No Motion"
$ws.Cells.Item(83, 6).Value = "Inactive"

$ws.Cells.Item(84, 1).Value = "2026-02-04"
$ws.Cells.Item(84, 2).Value = "14:07:54"
$ws.Cells.Item(84, 3).Value = "14:00"
$ws.Cells.Item(84, 4).Value = "Bathroom"
$ws.Cells.Item(84, 5).Value = "No Motion"
$ws.Cells.Item(84, 6).Value = "Inactive"

$ws.Cells.Item(85, 1).Value = "2026-02-04"
$ws.Cells.Item(85, 2).Value = "14:07:55"
$ws.Cells.Item(85, 3).Value = "14:00"
$ws.Cells.Item(85, 4).Value = "Bathroom"
$ws.Cells.Item(85, 5).Value = "No Motion"
$ws.Cells.Item(85, 6).Value = "Inactive"

$ws.Cells.Item(86, 1).Value = "2026-02-04"
$ws.Cells.Item(86, 2).Value = "14:07:56"
$ws.Cells.Item(86, 3).Value = "14:00"
$ws.Cells.Item(86, 4).Value = "Bathroom"
$ws.Cells.Item(86, 5).Value = "Motion Detected"
$ws.Cells.Item(86, 6).Value = "Active"

$ws.Cells.Item(87, 1).Value = "2026-02-04"
$ws.Cells.Item(87, 2).Value = "14:08:04"
$ws.Cells.Item(87, 3).Value = "14:00"
$ws.Cells.Item(87, 4).Value = "Bathroom"
$ws.Cells.Item(87, 5).Value = "No Motion"
$ws.Cells.Item(87, 6).Value = "Inactive"

$ws.Cells.Item(88, 1).Value = "2026-02-04"
$ws.Cells.Item(88, 2).Value = "14:08:09"
$ws.Cells.Item(88, 3).Value = "14:00"
$ws.Cells.Item(88, 4).Value = "Bathroom"
$ws.Cells.Item(88, 5).Value = "No Motion"
$ws.Cells.Item(88, 6).Value = "Inactive"

$ws.Cells.Item(89, 1).Value = "2026-02-04"
$ws.Cells.Item(89, 2).Value = "14:08:13"
$ws.Cells.Item(89, 3).Value = "14:00"
$ws.Cells.Item(89, 4).Value = "Bathroom"
$ws.Cells.Item(89, 5).Value = "Motion Detected"
$ws.Cells.Item(89, 6).Value = "Active"

$ws.Cells.Item(90, 1).Value = "2026-02-04"
$ws.Cells.Item(90, 2).Value = "14:08:21"
$ws.Cells.Item(90, 3).Value = "14:00"
$ws.Cells.Item(90, 4).Value = "Bathroom"
$ws.Cells.Item(90, 5).Value = "No Motion"
$ws.Cells.Item(90, 6).Value = "Inactive"

$ws.Cells.Item(91, 1).Value = "2026-02-04"
$ws.Cells.Item(91, 2).Value = "14:08:26"
$ws.Cells.Item(91, 3).Value = "14:00"
$ws.Cells.Item(91, 4).Value = "Bathroom"
$ws.Cells.Item(91, 5).Value = "Motion Detected"
$ws.Cells.Item(91, 6).Value = "Active"

$ws.Cells.Item(92, 1).Value = "2026-02-04"
$ws.Cells.Item(92, 2).Value = "14:08:34"
$ws.Cells.Item(92, 3).Value = "14:00"
$ws.Cells.Item(92, 4).Value = "Bathroom"
$ws.Cells.Item(92, 5).Value = "No Motion"
$ws.Cells.Item(92, 6).Value = "Inactive"

$ws.Cells.Item(93, 1).Value = "2026-02-04"
$ws.Cells.Item(93, 2).Value = "14:08:39"
$ws.Cells.Item(93, 3).Value = "14:00"
$ws.Cells.Item(93, 4).Value = "Bathroom"
$ws.Cells.Item(93, 5).Value = "No Motion"
$ws.Cells.Item(93, 6).Value = "Inactive"

$ws.Cells.Item(94, 1).Value = "2026-02-04"
$ws.Cells.Item(94, 2).Value = "14:08:44"
$ws.Cells.Item(94, 3).Value = "14:00"
$ws.Cells.Item(94, 4).Value = "Bathroom"
$ws.Cells.Item(94, 5).Value = "No Motion"
$ws.Cells.Item(94, 6).Value = "Inactive"

$ws.Cells.Item(95, 1).Value = "2026-02-04"
$ws.Cells.Item(95, 2).Value = "14:08:49"
$ws.Cells.Item(95, 3).Value = "14:00"
$ws.Cells.Item(95, 4).Value = "Bathroom"
$ws.Cells.Item(95, 5).Value = "No Motion"
$ws.Cells.Item(95, 6).Value = "Inactive"

$ws.Cells.Item(96, 1).Value = "2026-02-04"
$ws.Cells.Item(96, 2).Value = "14:08:49"
$ws.Cells.Item(96, 3).Value = "14:00"
$ws.Cells.Item(96, 4).Value = "Bathroom"
$ws.Cells.Item(96, 5).Value = "Motion Detected"
$ws.Cells.Item(96, 6).Value = "Active"


$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A63:A74").NumberFormat = "@"
$ws.Range("E63:E74").NumberFormat = "@"

$ws.Cells.Item(63, 1).Value = "2026-02-04"
$ws.Cells.Item(63, 2).Value = "14:07:51"
$ws.Cells.Item(63, 3).Value = "14:00"
$ws.Cells.Item(63, 4).Value = "Bathroom"
$ws.Cells.Item(63, 5).Value = "76.6%"
$ws.Cells.Item(63, 6).Value = "Active"

$ws.Cells.Item(64, 1).Value = "2026-02-04"
$ws.Cells.Item(64, 2).Value = "14:07:51"
$ws.Cells.Item(64, 3).Value = "14:00"
$ws.Cells.Item(64, 4).Value = "Bathroom"
$ws.Cells.Item(64, 5).Value = "77.5%"
$ws.Cells.Item(64, 6).Value = "Active"

$ws.Cells.Item(65, 1).Value = "2026-02-04"
$ws.Cells.Item(65, 2).Value = "14:07:52"
$ws.Cells.Item(65, 3).Value = "14:00"
$ws.Cells.Item(65, 4).Value = "Bathroom"
$ws.Cells.Item(65, 5).Value = "76.6%"
$ws.Cells.Item(65, 6).Value = "Active"

$ws.Cells.Item(66, 1).Value = "2026-02-04"
$ws.Cells.Item(66, 2).Value = "14:07:53"
$ws.Cells.Item(66, 3).Value = "14:00"
$ws.Cells.Item(66, 4).Value = "Bathroom"
$ws.Cells.Item(66, 5).Value = "77.4%"
$ws.Cells.Item(66, 6).Value = "Active"

$ws.Cells.Item(67, 1).Value = "2026-02-04"
$ws.Cells.Item(67, 2).Value = "14:07:54"
$ws.Cells.Item(67, 3).Value = "14:00"
$ws.Cells.Item(67, 4).Value = "Bathroom"
$ws.Cells.Item(67, 5).Value = "77.5%"
$ws.Cells.Item(67, 6).Value = "Active"

$ws.Cells.Item(68, 1).Value = "2026-02-04"
$ws.Cells.Item(68, 2).Value = "14:08:04"
$ws.Cells.Item(68, 3).Value = "14:00"
$ws.Cells.Item(68, 4).Value = "Bathroom"
$ws.Cells.Item(68, 5).Value = "77.5%"
$ws.Cells.Item(68, 6).Value = "Active"

$ws.Cells.Item(69, 1).Value = "2026-02-04"
$ws.Cells.Item(69, 2).Value = "14:08:09"
$ws.Cells.Item(69, 3).Value = "14:00"
$ws.Cells.Item(69, 4).Value = "Bathroom"
$ws.Cells.Item(69, 5).Value = "76.5%"
$ws.Cells.Item(69, 6).Value = "Active"

$ws.Cells.Item(70, 1).Value = "2026-02-04"
$ws.Cells.Item(70, 2).Value = "14:08:14"
$ws.Cells.Item(70, 3).Value = "14:00"
$ws.Cells.Item(70, 4).Value = "Bathroom"
$ws.Cells.Item(70, 5).Value = "77.5%"
$ws.Cells.Item(70, 6).Value = "Active"

$ws.Cells.Item(71, 1).Value = "2026-02-04"
$ws.Cells.Item(71, 2).Value = "14:08:19"
$ws.Cells.Item(71, 3).Value = "14:00"
$ws.Cells.Item(71, 4).Value = "Bathroom"
$ws.Cells.Item(71, 5).Value = "76.6%"
$ws.Cells.Item(71, 6).Value = "Active"

$ws.Cells.Item(72, 1).Value = "2026-02-04"
$ws.Cells.Item(72, 2).Value = "14:08:34"
$ws.Cells.Item(72, 3).Value = "14:00"
$ws.Cells.Item(72, 4).Value = "Bathroom"
$ws.Cells.Item(72, 5).Value = "77.5%"
$ws.Cells.Item(72, 6).Value = "Active"

$ws.Cells.Item(73, 1).Value = "2026-02-04"
$ws.Cells.Item(73, 2).Value = "14:08:39"
$ws.Cells.Item(73, 3).Value = "14:00"
$ws.Cells.Item(73, 4).Value = "Bathroom"
$ws.Cells.Item(73, 5).Value = "76.6%"
$ws.Cells.Item(73, 6).Value = "Active"

$ws.Cells.Item(74, 1).Value = "2026-02-04"
$ws.Cells.Item(74, 2).Value = "14:08:44"
$ws.Cells.Item(74, 3).Value = "14:00"
$ws.Cells.Item(74, 4).Value = "Bathroom"
$ws.Cells.Item(74, 5).Value = "77.5%"
$ws.Cells.Item(74, 6).Value = "Active"


$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A63:A74").NumberFormat = "@"

$ws.Cells.Item(63, 1).Value = "2026-02-04"
$ws.Cells.Item(63, 2).Value = "14:07:51"
$ws.Cells.Item(63, 3).Value = "14:00"
$ws.Cells.Item(63, 4).Value = "Bathroom"
$ws.Cells.Item(63, 5).Value = "24.8C"
$ws.Cells.Item(63, 6).Value = "Active"

$ws.Cells.Item(64, 1).Value = "2026-02-04"
$ws.Cells.Item(64, 2).Value = "14:07:52"
$ws.Cells.Item(64, 3).Value = "14:00"
$ws.Cells.Item(64, 4).Value = "Bathroom"
$ws.Cells.Item(64, 5).Value = "24.8C"
$ws.Cells.Item(64, 6).Value = "Active"

$ws.Cells.Item(65, 1).Value = "2026-02-04"
$ws.Cells.Item(65, 2).Value = "14:07:52"
$ws.Cells.Item(65, 3).Value = "14:00"
$ws.Cells.Item(65, 4).Value = "Bathroom"
$ws.Cells.Item(65, 5).Value = "24.8C"
$ws.Cells.Item(65, 6).Value = "Active"

$ws.Cells.Item(66, 1).Value = "2026-02-04"
$ws.Cells.Item(66, 2).Value = "14:07:53"
$ws.Cells.Item(66, 3).Value = "14:00"
$ws.Cells.Item(66, 4).Value = "Bathroom"
$ws.Cells.Item(66, 5).Value = "24.8C"
$ws.Cells.Item(66, 6).Value = "Active"

$ws.Cells.Item(67, 1).Value = "2026-02-04"
$ws.Cells.Item(67, 2).Value = "14:07:54"
$ws.Cells.Item(67, 3).Value = "14:00"
$ws.Cells.Item(67, 4).Value = "Bathroom"
$ws.Cells.Item(67, 5).Value = "24.8C"
$ws.Cells.Item(67, 6).Value = "Active"

$ws.Cells.Item(68, 1).Value = "2026-02-04"
$ws.Cells.Item(68, 2).Value = "14:08:04"
$ws.Cells.Item(68, 3).Value = "14:00"
$ws.Cells.Item(68, 4).Value = "Bathroom"
$ws.Cells.Item(68, 5).Value = "24.8C"
$ws.Cells.Item(68, 6).Value = "Active"

$ws.Cells.Item(69, 1).Value = "2026-02-04"
$ws.Cells.Item(69, 2).Value = "14:08:09"
$ws.Cells.Item(69, 3).Value = "14:00"
$ws.Cells.Item(69, 4).Value = "Bathroom"
$ws.Cells.Item(69, 5).Value = "24.8C"
$ws.Cells.Item(69, 6).Value = "Active"

$ws.Cells.Item(70, 1).Value = "2026-02-04"
$ws.Cells.Item(70, 2).Value = "14:08:14"
$ws.Cells.Item(70, 3).Value = "14:00"
$ws.Cells.Item(70, 4).Value = "Bathroom"
$ws.Cells.Item(70, 5).Value = "24.8C"
$ws.Cells.Item(70, 6).Value = "Active"

$ws.Cells.Item(71, 1).Value = "2026-02-04"
$ws.Cells.Item(71, 2).Value = "14:08:20"
$ws.Cells.Item(71, 3).Value = "14:00"
$ws.Cells.Item(71, 4).Value = "Bathroom"
$ws.Cells.Item(71, 5).Value = "24.8C"
$ws.Cells.Item(71, 6).Value = "Active"

$ws.Cells.Item(72, 1).Value = "2026-02-04"
$ws.Cells.Item(72, 2).Value = "14:08:34"
$ws.Cells.Item(72, 3).Value = "14:00"
$ws.Cells.Item(72, 4).Value = "Bathroom"
$ws.Cells.Item(72, 5).Value = "24.8C"
$ws.Cells.Item(72, 6).Value = "Active"

$ws.Cells.Item(73, 1).Value = "2026-02-04"
$ws.Cells.Item(73, 2).Value = "14:08:40"
$ws.Cells.Item(73, 3).Value = "14:00"
$ws.Cells.Item(73, 4).Value = "Bathroom"
$ws.Cells.Item(73, 5).Value = "24.8C"
$ws.Cells.Item(73, 6).Value = "Active"

$ws.Cells.Item(74, 1).Value = "2026-02-04"
$ws.Cells.Item(74, 2).Value = "14:08:45"
$ws.Cells.Item(74, 3).Value = "14:00"
$ws.Cells.Item(74, 4).Value = "Bathroom"
$ws.Cells.Item(74, 5).Value = "24.8C"
$ws.Cells.Item(74, 6).Value = "Active"

